$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D sometimes hold plain numeric-looking text (e.g. "269.92").
# Excel auto-converts such strings typed via .Value into real numbers, so for
# those we force a Text number format, assign the value, then restore the
# default "Normal" style so the cell keeps looking like the untouched cells
# around it once the edit is done.

$ws.Range('D2').Value = '43.863.78'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '2.295.88'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '113.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +17.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '269.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.66%  '
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.618'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '48.13'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0947'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.04'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +14.57%  '
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.85'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.65%  '
$ws.Range('D15').Value = '2.640.09'
$ws.Range('E15').Value = '  +0.33%  '
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').Value = '2.297.86'
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('D18').Value = '43.772.50'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000110'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.86'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +11.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('E22').Value = '  -2.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.03'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +11.55%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '233.20'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.68'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.04%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.61'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '41.74'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +8.13%  '
$ws.Range('E29').Value = '  -1.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.27'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '175.53'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('E32').Value = '  +4.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.53'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.71'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.44%  '
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.65'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.74%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0364'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.76%  '
$ws.Range('E38').Value = '  +0.53%  '
$ws.Range('E39').Value = '  +7.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '75.11'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +16.68%  '
$ws.Range('E41').Value = '  +3.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.84'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +13.02%  '
$ws.Range('E43').Value = '  +2.78%  '
$ws.Range('E44').Value = '  +22.19%  '
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.39'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.84'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0997'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '101.52'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.64%  '
$ws.Range('E50').Value = '  +3.35%  '
$ws.Range('E51').Value = '  +5.96%  '
